# Applies the updated crypto price/volume figures captured in the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / percentage cells: safe to assign directly ---
$directUpdates = @{
    'D2' = '38.698.11'
    'E2' = '  -5.24%  '
    'D3' = '2.207.43'
    'E3' = '  -7.43%  '
    'E4' = '  +0.15%  '
    'E5' = '  -5.49%  '
    'E6' = '  -10.31%  '
    'E7' = '  -5.16%  '
    'E8' = '  +0.10%  '
    'E9' = '  -7.69%  '
    'E10' = '  -6.90%  '
    'E11' = '  -11.26%  '
    'E12' = '  -13.62%  '
    'E13' = '  -1.89%  '
    'D14' = '2.551.45'
    'E14' = '  -7.25%  '
    'E15' = '  -7.73%  '
    'E16' = '  -8.46%  '
    'D17' = '2.226.30'
    'E17' = '  -6.63%  '
    'E18' = '  -7.32%  '
    'D19' = '38.655.88'
    'E19' = '  -5.10%  '
    'D20' = '0.0₃0854'
    'E20' = '  -6.57%  '
    'E21' = '  -8.22%  '
    'E22' = '  -7.03%  '
    'E23' = '  -10.05%  '
    'E24' = '  -4.43%  '
    'E25' = '  -0.11%  '
    'E26' = '  -10.42%  '
    'E27' = '  -6.65%  '
    'E28' = '  -1.23%  '
    'E29' = '  -7.52%  '
    'E30' = '  -6.16%  '
    'E31' = '  -4.65%  '
    'E32' = '  -8.87%  '
    'E33' = '  -0.02%  '
    'E34' = '  -9.18%  '
    'E35' = '  -5.04%  '
    'E36' = '  -7.20%  '
    'E37' = '  -4.95%  '
    'E38' = '  -6.74%  '
    'E39' = '  -5.42%  '
    'E40' = '  -8.80%  '
    'E41' = '  -11.62%  '
    'E42' = '  -6.53%  '
    'D43' = '1.895.87'
    'E43' = '  -3.00%  '
    'E44' = '  -12.04%  '
    'E45' = '  -7.13%  '
    'B46' = 'EnergySwap'
    'C46' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E46' = '  -8.33%  '
    'B47' = 'FraxShare'
    'C47' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E47' = '  -5.02%  '
    'E48' = '  -11.06%  '
    'D49' = '2.431.64'
    'E49' = '  -7.09%  '
    'E50' = '  -7.28%  '
    'E51' = '  -7.57%  '
}
foreach ($addr in $directUpdates.Keys) {
    $ws.Range($addr).Value = $directUpdates[$addr]
}

# --- "Price" cells whose new text is a plain number (e.g. "79.34"): assigning
# such a string straight to .Value would make Excel auto-convert it to a real
# number, losing the original text formatting. Route them through a helper
# cell formatted as Text and paste-special the value across so the destination
# keeps its original (default) style, matching how the source text cells are.
$numericTextUpdates = @{
    'D5' = '295.94'
    'D6' = '79.34'
    'D7' = '0.503'
    'D9' = '0.456'
    'D10' = '0.0769'
    'D11' = '27.73'
    'D12' = '45.91'
    'D15' = '6.08'
    'D16' = '13.89'
    'D18' = '0.709'
    'D21' = '5.70'
    'D22' = '64.56'
    'D23' = '9.76'
    'D24' = '223.71'
    'D27' = '1.70'
    'D29' = '22.02'
    'D30' = '8.82'
    'D31' = '148.29'
    'D32' = '30.95'
    'D34' = '4.74'
    'D39' = '0.0945'
    'D41' = '14.28'
    'D44' = '2.00'
    'D46' = '16.08'
    'D47' = '8.92'
    'D50' = '67.83'
    'D51' = '86.81'
}
$helper = $ws.Range("Z1")
foreach ($addr in $numericTextUpdates.Keys) {
    $helper.NumberFormat = "@"
    $helper.Value = $numericTextUpdates[$addr]
    $helper.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}
$helper.Clear()
$excel.CutCopyMode = $false
